$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 679.4
$ws.Range("C2").Value = 671.45
$ws.Range("D2").Value = 675.5
$ws.Range("E2").Value = 675.75
$ws.Range("G2").Value = 673.95
$ws.Range("B3").Value = 3275
$ws.Range("D3").Value = 3269
$ws.Range("E3").Value = 3261.75
$ws.Range("F3").Value = 42
$ws.Range("G3").Value = 3227.8
$ws.Range("B4").Value = 478
$ws.Range("C4").Value = 472.1
$ws.Range("D4").Value = 476.6
$ws.Range("E4").Value = 476.45
$ws.Range("F4").Value = 17
$ws.Range("G4").Value = 474.4
$ws.Range("B5").Value = 1596
$ws.Range("C5").Value = 1575.35
$ws.Range("D5").Value = 1587.05
$ws.Range("E5").Value = 1591.75
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 1582.6
$ws.Range("B6").Value = 7370
$ws.Range("C6").Value = 7247
$ws.Range("D6").Value = 7340.5
$ws.Range("E6").Value = 7341.55
$ws.Range("F6").Value = 11
$ws.Range("G6").Value = 7300
$ws.Range("B7").Value = 195.9
$ws.Range("C7").Value = 193.53
$ws.Range("D7").Value = 194.15
$ws.Range("E7").Value = 194.34
$ws.Range("F7").Value = 51
$ws.Range("G7").Value = 195.6
$ws.Range("B8").Value = 287.5
$ws.Range("C8").Value = 281.75
$ws.Range("D8").Value = 286
$ws.Range("E8").Value = 286.25
$ws.Range("F8").Value = 139
$ws.Range("G8").Value = 282.35
$ws.Range("B9").Value = 492.4
$ws.Range("C9").Value = 485.7
$ws.Range("D9").Value = 487.5
$ws.Range("E9").Value = 486.95
$ws.Range("F9").Value = 59
$ws.Range("G9").Value = 487.55
$ws.Range("B10").Value = 883.65
$ws.Range("C10").Value = 871.9
$ws.Range("D10").Value = 877.25
$ws.Range("E10").Value = 878.6
$ws.Range("F10").Value = 30
$ws.Range("G10").Value = 872.85
$ws.Range("B11").Value = 4945
$ws.Range("C11").Value = 4801.5
$ws.Range("D11").Value = 4929
$ws.Range("E11").Value = 4935.1
$ws.Range("G11").Value = 4836
$ws.Range("B12").Value = 174.89
$ws.Range("C12").Value = 172.51
$ws.Range("D12").Value = 174.15
$ws.Range("E12").Value = 174.4
$ws.Range("F12").Value = 111
$ws.Range("G12").Value = 173.23
$ws.Range("B13").Value = 1439.75
$ws.Range("C13").Value = 1424.4
$ws.Range("D13").Value = 1431.2
$ws.Range("E13").Value = 1431.05
$ws.Range("F13").Value = 30
$ws.Range("G13").Value = 1435.55
$ws.Range("B14").Value = 1599
$ws.Range("C14").Value = 1579.05
$ws.Range("D14").Value = 1595.6
$ws.Range("E14").Value = 1596.9
$ws.Range("F14").Value = 127
$ws.Range("G14").Value = 1582.25
$ws.Range("B15").Value = 686.5
$ws.Range("C15").Value = 678.9
$ws.Range("D15").Value = 682.1
$ws.Range("E15").Value = 683.6
$ws.Range("F15").Value = 83
$ws.Range("G15").Value = 679.9
$ws.Range("B16").Value = 1112.95
$ws.Range("C16").Value = 1099.45
$ws.Range("D16").Value = 1105.05
$ws.Range("E16").Value = 1105.65
$ws.Range("F16").Value = 113
$ws.Range("G16").Value = 1101
$ws.Range("B17").Value = 1515.9
$ws.Range("C17").Value = 1495.95
$ws.Range("D17").Value = 1498.5
$ws.Range("E17").Value = 1502.35
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 1498.4
$ws.Range("B18").Value = 1495.2
$ws.Range("C18").Value = 1486.05
$ws.Range("D18").Value = 1488.5
$ws.Range("E18").Value = 1488.9
$ws.Range("F18").Value = 42
$ws.Range("G18").Value = 1489.6
$ws.Range("B19").Value = 1059.25
$ws.Range("C19").Value = 1036
$ws.Range("D19").Value = 1052
$ws.Range("E19").Value = 1052.45
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 1038.2
$ws.Range("B20").Value = 734.7
$ws.Range("C20").Value = 717.2
$ws.Range("D20").Value = 734
$ws.Range("E20").Value = 731.65
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 718.55
$ws.Range("B21").Value = 2946
$ws.Range("C21").Value = 2865
$ws.Range("D21").Value = 2924
$ws.Range("E21").Value = 2928.6
$ws.Range("F21").Value = 35
$ws.Range("G21").Value = 2872.35
$ws.Range("B22").Value = 300
$ws.Range("C22").Value = 293.2
$ws.Range("D22").Value = 298.95
$ws.Range("E22").Value = 298.95
$ws.Range("F22").Value = 46
$ws.Range("G22").Value = 293.85
$ws.Range("B23").Value = 369.45
$ws.Range("C23").Value = 366.1
$ws.Range("D23").Value = 368.5
$ws.Range("E23").Value = 368.45
$ws.Range("F23").Value = 102
$ws.Range("G23").Value = 366.5
$ws.Range("B24").Value = 2959.35
$ws.Range("C24").Value = 2914.45
$ws.Range("D24").Value = 2951.45
$ws.Range("E24").Value = 2955.1
$ws.Range("F24").Value = 40
$ws.Range("G24").Value = 2917.2
$ws.Range("B25").Value = 845.8
$ws.Range("C25").Value = 835.1
$ws.Range("D25").Value = 839.2
$ws.Range("E25").Value = 839.2
$ws.Range("F25").Value = 129
$ws.Range("G25").Value = 842.85
$ws.Range("B26").Value = 767.05
$ws.Range("C26").Value = 751
$ws.Range("D26").Value = 760
$ws.Range("E26").Value = 757.55
$ws.Range("F26").Value = 10
$ws.Range("G26").Value = 760.2
$ws.Range("B27").Value = 1118.75
$ws.Range("C27").Value = 1100.7
$ws.Range("D27").Value = 1105
$ws.Range("E27").Value = 1105.85
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 1112.05
$ws.Range("B28").Value = 997.25
$ws.Range("C28").Value = 981.4
$ws.Range("D28").Value = 992
$ws.Range("E28").Value = 993.4
$ws.Range("F28").Value = 115
$ws.Range("G28").Value = 982.3
$ws.Range("B29").Value = 453.6
$ws.Range("C29").Value = 447.5
$ws.Range("D29").Value = 448.5
$ws.Range("E29").Value = 448.65
$ws.Range("F29").Value = 78
$ws.Range("G29").Value = 449.1
$ws.Range("B30").Value = 183.5
$ws.Range("C30").Value = 181.4
$ws.Range("D30").Value = 182.9
$ws.Range("E30").Value = 183.15
$ws.Range("F30").Value = 285
$ws.Range("G30").Value = 181.75
$ws.Range("B31").Value = 11271
$ws.Range("C31").Value = 11167.6
$ws.Range("D31").Value = 11250
$ws.Range("E31").Value = 11242.8
$ws.Range("G31").Value = 11175
